$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.969.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.622.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.77"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.620.89"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000188"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.104.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "71.860.98"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.691.59"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.12"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.85"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.21"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.773.95"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0942"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.32"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "483.80"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.31"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.20"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.90"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.71"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.325"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.04"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.41"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.63"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.538"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.66"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.603"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.08%  "
